$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the updated coin-ranking snapshot values.
# For cells whose new text looks like a plain number or a percentage,
# the cell's NumberFormat is first set to Text ("@") so Excel keeps
# the value as a literal string (matching the source inlineStr cells)
# instead of silently re-interpreting it as a numeric/percentage value.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '308.38'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '0.46%'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '40.80'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '1.57%'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '5.117'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '-0.68%'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.07613'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '-1.12%'
$ws.Range('B6').Value = 'GateToken'
$ws.Range('C6').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '4.254'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '-0.27%'
$ws.Range('B7').Value = 'FTXToken'
$ws.Range('C7').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.623'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '-0.19%'
$ws.Range('B8').Value = 'BTSEToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '2.448'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '0.82%'
$ws.Range('B9').Value = 'MXToken'
$ws.Range('C9').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.9002'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '2.49%'
$ws.Range('B10').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C10').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.1084'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '8.86%'
$ws.Range('B11').Value = 'WazirX'
$ws.Range('C11').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.1760'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '1.54%'
$ws.Range('B12').Value = 'MandalaExchangeToken'
$ws.Range('C12').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.09195'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '3.39%'
$ws.Range('B13').Value = 'BitrueCoin'
$ws.Range('C13').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.04185'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '-4.98%'
$ws.Range('B14').Value = 'BitMartToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.1052'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '-0.49%'
$ws.Range('B15').Value = 'BitForexToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.001252'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '-0.69%'
$ws.Range('B16').Value = 'TigerCash'
$ws.Range('C16').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.005896'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '-0.47%'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.353'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '-0.08%'
$ws.Range('B18').Value = 'BitpandaEcosystemToken'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.3294'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '-0.21%'
$ws.Range('B19').Value = 'MCDex'
$ws.Range('C19').Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.550'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '-6.16%'
$ws.Range('B20').Value = 'ProBitToken'
$ws.Range('C20').Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.1365'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '1.91%'
$ws.Range('B21').Value = 'ZBToken'
$ws.Range('C21').Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.2682'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '-14.48%'
$ws.Range('B22').Value = 'CoinExToken'
$ws.Range('C22').Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.04088'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '-1.38%'
$ws.Range('B23').Value = 'BitKan'
$ws.Range('C23').Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.001223'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '2.07%'
$ws.Range('B24').Value = 'HotbitToken'
$ws.Range('C24').Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.004087'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '0.62%'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '6.70%'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02366'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '1.04%'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.05181'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '0.79%'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.007776'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '-2.05%'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.1300'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '-1.81%'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.006786'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '6.77%'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '0.52%'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.008556'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '-0.46%'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '0.66%'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.00006945'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '6.43%'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.00000000751'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '0.07%'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.02897'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '753.92%'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.004204'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '-39.96%'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.00002102'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '0.07%'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0002002'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '0.07%'
